$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "attr4" header to "m_mymap"
$ws.Range("F1").Value = "m_mymap"

# Add the new "t_test" column with a "hello" value for the first data row
$ws.Range("M1").Value = "t_test"
$ws.Range("M2").Value = "hello"

# Match the final cursor/selection position recorded in the saved file
$ws.Range("J9").Select() | Out-Null
